## Modify: SystemVerilog_語法學習.pptx
## 1. Add Function Coverage slide

$p = $ppt.ActivePresentation

# --- 1. Add a new bullet "Function Coverage" to the 目錄 (table of contents)
#        slide, right after the existing "Communication" bullet.
$tocSlide = $p.Slides.Item(2)
$tocBody  = $tocSlide.Shapes.Item(2).TextFrame.TextRange
$commParagraph = $tocBody.Paragraphs(4)
$commParagraph.InsertAfter("`rFunction Coverage") | Out-Null

# --- 2. Add a new "Chapter 5 / Function Coverage" section-header slide at the
#        end of the deck (same pattern as the existing Chapter 1-4 slides).
#        Duplicate the most similar existing chapter slide (Chapter 3 /
#        Processes, slide 33) so the new slide inherits the correct layout,
#        placeholders and formatting, then retarget its text.
$chapterTemplate = $p.Slides.Item(33)
$chapterTemplate.Copy()
$p.Slides.Paste($p.Slides.Count + 1) | Out-Null
$newChapterSlide = $p.Slides.Item($p.Slides.Count)
$newChapterSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Chapter 5"
$newChapterSlide.Shapes.Item(2).TextFrame.TextRange.Text = "Function Coverage"
